$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 154 (pushes the existing rows 154..253 down to 155..254,
# matching the diff where every subsequent row's data shifts down by one position
# and a brand-new record is inserted at the top of that block).
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new record's data.
$ws.Cells.Item(154, 1).Value  = 3
$ws.Cells.Item(154, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(154, 3).Value  = "Coquimbo"
$ws.Cells.Item(154, 4).Value  = 45236
$ws.Cells.Item(154, 5).Value  = 5
$ws.Cells.Item(154, 6).Value  = 100112010
$ws.Cells.Item(154, 7).Value  = "Achicoria"
$ws.Cells.Item(154, 8).Value  = "Sin especificar"
$ws.Cells.Item(154, 9).Value  = "Primera"
$ws.Cells.Item(154, 10).Value = 60
$ws.Cells.Item(154, 11).Value = 7000
$ws.Cells.Item(154, 12).Value = 7000
$ws.Cells.Item(154, 13).Value = 7000
$ws.Cells.Item(154, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(154, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(154, 16).Value = 438
$ws.Cells.Item(154, 17).Value = 16
$ws.Cells.Item(154, 18).Value = "Hortaliza"
